$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column B (abg_hypercap_threshold) values
$ws.Range("B2").Value = 8.84
$ws.Range("B4").Value = 8.77
$ws.Range("B5").Value = 9.140000000000001
$ws.Range("B6").Value = 11.71
$ws.Range("B7").Value = 1.6
$ws.Range("B8").Value = 4.84
$ws.Range("B9").Value = 5.83
$ws.Range("B10").Value = 14.84
$ws.Range("B11").Value = 22.18
$ws.Range("B12").Value = 2.65

# Update column E (vbg_hypercap_threshold) values
$ws.Range("E3").Value = 6.27
$ws.Range("E4").Value = 6.41
$ws.Range("E5").Value = 9.630000000000001
$ws.Range("E6").Value = 14.89
$ws.Range("E7").Value = 3.17
$ws.Range("E8").Value = 6.3
$ws.Range("E9").Value = 5.67
$ws.Range("E10").Value = 14.02
$ws.Range("E11").Value = 25.41
$ws.Range("E12").Value = 2.17
